$d = $word.ActiveDocument

$replacements = @(
    @("36×43=", "15×35="),
    @("67×11=", "22×93="),
    @("16×29=", "93×64="),
    @("71×63=", "60×51="),
    @("38×77=", "52×73="),
    @("58×96=", "15×41="),
    @("71×51=", "61×50="),
    @("27×93=", "49×95="),
    @("46×52=", "72×88="),
    @("81×64=", "86×83="),
    @("62×78=", "87×26="),
    @("44×60=", "89×20="),
    @("45×20=", "29×93="),
    @("59×82=", "11×19="),
    @("49×14=", "42×25="),
    @("13×72=", "95×14="),
    @("81×73=", "67×37="),
    @("85×26=", "37×59="),
    @("50×19=", "97×75="),
    @("19×56=", "71×97="),
    @("95×37=", "92×63="),
    @("42×83=", "42×22="),
    @("86×24=", "74×91="),
    @("73×57=", "28×67="),
    @("94×57=", "47×31=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
